# Update the "Library preparation protocol" sheet to match the updated
# subscription query for 10x v2 data:
#   - Insert a new "nucleic_acid_source" column (library_construction_approach
#     group) with value "single cell"
#   - Change the existing "End bias" value from "full length" to "3 prime tag"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library preparation protocol")

# Column M (13) currently holds library_construction_approach.text; insert a
# new blank column there so the new nucleic_acid_source field lives at M,
# pushing library_construction_approach (and everything after it) one column
# to the right.
$ws.Columns.Item(13).Insert()

# Row 1: long-form field description (wraps the header row)
$ws.Cells.Item(1, 13).Value = "Source cells or organelles from which nucleic acid molecules were collected."

# Row 2 & Row 4: machine-readable field name (row 4 mirrors row 2 in this sheet)
$ws.Cells.Item(2, 13).Value = "library_preparation_protocol_json.nucleic_acid_source"
$ws.Cells.Item(4, 13).Value = "library_preparation_protocol_json.nucleic_acid_source"

# Row 6: example/default data value for the new field
$ws.Cells.Item(6, 13).Value = 'single cell"'

# The "End bias" example value (previously column AB, now shifted to AC after
# the column insert above) changes from "full length" to "3 prime tag".
$ws.Cells.Item(6, 29).Value = "3 prime tag"
